# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" right after "总计", pushing
#   "2021-Q1" and "2020-Q4" one tab to the right.
# - Populate "总计" with a new row for 2022-Q4 (and re-derive the
#   existing rows' running index), shifting 2021-Q1/2020-Q4 down.
# - Populate the new "2022-Q4" sheet with the fund holdings for that
#   quarter (same layout as the other quarter sheets, but the D1
#   header reads "基金规模" instead of "基金金额").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. 总计 (summary) sheet: insert the 2022-Q4 row at the top of the
#    data and push the older rows down by one.
# ---------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Push 2021-Q1 / 2020-Q4 rows down (row3 -> row4, row2 -> row3),
# carrying formatting along via Copy so the styled A-column cell
# (index style) follows the data.
$wsTotal.Range("A3:D3").Copy($wsTotal.Range("A4:D4"))
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.05

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q1"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2020-Q4"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.01

# ---------------------------------------------------------------
# 2. Add the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($null, $wsTotal)
$wsNew.Name = "2022-Q4"

# Match the page margins used by the sibling quarter sheets.
$wsNew.PageSetup.LeftMargin = 54
$wsNew.PageSetup.RightMargin = 54
$wsNew.PageSetup.TopMargin = 72
$wsNew.PageSetup.BottomMargin = 72
$wsNew.PageSetup.HeaderMargin = 36
$wsNew.PageSetup.FooterMargin = 36

# Clone the bold/boxed header + index-column styling from the
# "2021-Q1" sheet, which is laid out exactly the same way.
$ws2021 = $wb.Worksheets.Item("2021-Q1")
$ws2021.Range("B1:H1").Copy($wsNew.Range("B1:H1"))
$ws2021.Range("A2:A3").Copy($wsNew.Range("A2:A3"))

# Header row.
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Numeric-looking text columns need a text format so values like
# "004945" / "2.08" / "92.79" are stored verbatim (not coerced to
# numbers, which would drop the leading/trailing zeros). Apply the
# format per contiguous block - a multi-area Range here doesn't
# propagate NumberFormat reliably to every member cell.
$wsNew.Range("B2:B3").NumberFormat = "@"
$wsNew.Range("D2:G3").NumberFormat = "@"

$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "004945"
$wsNew.Range("C2").Value = "长信中证500指数增强A"
$wsNew.Range("D2").Value = "2.08"
$wsNew.Range("E2").Value = "92.79"
$wsNew.Range("F2").Value = "1.70"
$wsNew.Range("G2").Value = "0.0354"
$wsNew.Range("H2").Value = 2

$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "013881"
$wsNew.Range("C3").Value = "长信中证500指数增强C"
$wsNew.Range("D3").Value = "0.96"
$wsNew.Range("E3").Value = "92.79"
$wsNew.Range("F3").Value = "1.70"
$wsNew.Range("G3").Value = "0.0163"
$wsNew.Range("H3").Value = 2

# Drop the incidental "Text" style the NumberFormat="@" step left
# behind so these cells fall back to the default style, same as in
# the sibling quarter sheets.
$wsNew.Range("B2:B3").Style = "Normal"
$wsNew.Range("D2:G3").Style = "Normal"


